$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B:D hold numeric-looking odds/stat values that the source workbook
# stores as text (inline strings), not numbers - force Text format first so
# assigning e.g. "9.75" does not get auto-coerced into a Number cell.
$ws.Range("B2:D36").NumberFormat = "@"

$fixtures = @(
    @("VfL Bochum 1848 - SV Werder Bremen", "1.42", "4.25", "9.75", "Bundesliga"),
    @("Borussia Monchengladbach - Borussia Dortmund", "2.18", "3.08", "9.08", "Bundesliga"),
    @("VfL Wolfsburg - 1. FSV Mainz 05", "1.44", "3.92", "9.67", "Bundesliga"),
    @("TSG Hoffenheim - SC Freiburg", "1.58", "2.83", "8.25", "Bundesliga"),
    @("Crystal Palace - Manchester City", "1.60", "4.33", "9.83", "Premier League"),
    @("Aston Villa - Southampton", "0.93", "4.00", "11.17", "Premier League"),
    @("Leicester City - Brighton Hove Albion", "1.63", "3.42", "9.75", "Premier League"),
    @("Fulham - Arsenal", "1.39", "3.75", "10.92", "Premier League"),
    @("Tottenham Hotspur - Chelsea", "1.36", "4.17", "11.33", "Premier League"),
    @("West Ham United - Wolverhampton Wanderers", "1.35", "3.83", "9.58", "Premier League"),
    @("Sheffield Wednesday - Preston North End", "1.15", "3.08", "10.33", "EFL Championship"),
    @("Leeds United - Derby County", "1.71", "4.25", "9.08", "EFL Championship"),
    @("Luton Town - Swansea City", "1.46", "3.67", "8.83", "EFL Championship"),
    @("Cardiff City - Watford", "1.28", "4.17", "11.33", "EFL Championship"),
    @("Queens Park Rangers - Norwich City", "1.92", "4.08", "12.25", "EFL Championship"),
    @("Millwall - Coventry City", "1.30", "4.33", "10.25", "EFL Championship"),
    @("Hull City - Blackburn Rovers", "1.88", "3.42", "10.83", "EFL Championship"),
    @("Inter Milan - Parma", "1.54", "3.58", "9.75", "Serie A"),
    @("Atalanta - AC Milan", "1.56", "3.92", "8.08", "Serie A"),
    @("Genoa - Torino", "2.38", "4.92", "9.17", "Serie A"),
    @("Juventus - Bologna", "1.48", "3.92", "7.25", "Serie A"),
    @("AS Roma - Lecce", "1.17", "4.67", "8.92", "Serie A"),
    @("RC Celta - RCD Mallorca", "1.22", "4.33", "7.75", "La Liga"),
    @("Valencia CF - Rayo Vallecano", "1.38", "5.92", "10.75", "La Liga"),
    @("Girona FC - Real Madrid", "1.61", "4.67", "7.50", "La Liga"),
    @("Athletic Club - Villarreal CF", "1.56", "6.17", "8.50", "La Liga"),
    @("LOSC Lille - Stade Brestois 29", "0.98", "4.33", "8.83", "Ligue 1"),
    @("RC Strasbourg Alsace - Stade DE Reims", "1.35", "3.50", "8.42", "Ligue 1"),
    @("AS Saint-Étienne - Marseille", "1.79", "4.00", "9.75", "Ligue 1"),
    @("Moreirense - Sporting CP", "0.98", "5.17", "9.25", "Primera Division"),
    @("Sporting Braga - Estoril", "1.19", "4.25", "11.58", "Primera Division"),
    @("Benfica - Vitoria Guimaraes", "1.38", "4.75", "10.92", "Primera Division"),
    @("Santa Clara - Rio Ave", "1.31", "7.08", "9.25", "Primera Division"),
    @("FC Famalicao - FC Porto", "2.08", "5.25", "9.00", "Primera Division"),
    @("Boavista FC - SC Farense", "1.21", "4.50", "10.75", "Primera Division")
)

for ($i = 0; $i -lt $fixtures.Count; $i++) {
    $row = $i + 2
    $vals = $fixtures[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}

